$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (municipio-nombre) metadata updates
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column G (aragon) metadata updates
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("G4").Value = "URI-Comunidad"
$ws.Range("G5").Clear()
